$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number by Excel;
# force text format, assign, then restore the default "Normal" style so no
# spurious style index is left on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4027"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08223"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.106"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.406"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.277"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001096"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06497"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.933"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.195"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.280"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.127"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1037"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.005"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.778"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02438"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.331"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06424"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2160"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.885"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.199"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6436"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.217"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.189"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5992"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.634"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.212"
$ws.Range("D50").Style = "Normal"

# Remaining cells are safe to assign directly (they are not numeric-looking).
$ws.Range("D2").Value = "30.011.95"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.913.39"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  -2.94%  "
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.913.65"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").Value = "30.062.15"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").Value = "2.132.01"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("E29").Value = "  -4.67%  "
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  +3.90%  "
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("E36").Value = "  +3.06%  "
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E43").Value = "  -4.76%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("E46").Value = "  +5.96%  "
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("E50").Value = "  -2.41%  "
